# Rebuild the summary reference table headers for the organization
# identifier columns so they use the same friendly titles as the rest
# of the workbook instead of raw flatten-tool JSONPath placeholders.
#
# fun_location!B1:  "fundingOrganization/0/Identifier"   -> "Funding Org:Identifier"
# rec_location!B1:  "recipientOrganization/0/Identifier" -> "Recipient Org:Identifier"

$wb = $excel.ActiveWorkbook

$wsFunLocation = $wb.Worksheets.Item("fun_location")
$wsFunLocation.Range("B1").Value = "Funding Org:Identifier"

$wsRecLocation = $wb.Worksheets.Item("rec_location")
$wsRecLocation.Range("B1").Value = "Recipient Org:Identifier"
